$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("asistencia")
$ws.Range("A1").Value = "test"
